{"js": "// Update the date heading and every two-digit \u00f7 one-digit practice\n// problem in the worksheet table with the next day's regenerated set.\n\nconst body = context.document.body;\n\n// --- 1. Update the date paragraph above the table -----------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(\"2024-09-16 Monday\", Word.InsertLocation.replace);\n\n// --- 2. Update the practice-problem table --------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Each 4-row \"block\" in the table holds one row of 5 problems\n// followed by 3 blank rows. Replacements are listed row-by-row,\n// left-to-right, in document order.\nconst replacements = [\n  [\"23\u00f74=5, 3\", \"82\u00f78=10, 2\", \"89\u00f78=11, 1\", \"82\u00f77=11, 5\", \"27\u00f73=9, 0\"],\n  [\"71\u00f75=14, 1\", \"71\u00f78=8, 7\", \"38\u00f72=19, 0\", \"60\u00f77=8, 4\", \"91\u00f73=30, 1\"],\n  [\"25\u00f73=8, 1\", \"79\u00f74=19, 3\", \"28\u00f76=4, 4\", \"80\u00f72=40, 0\", \"92\u00f79=10, 2\"],\n  [\"16\u00f79=1, 7\", \"22\u00f78=2, 6\", \"62\u00f77=8, 6\", \"80\u00f77=11, 3\", \"60\u00f78=7, 4\"],\n  [\"94\u00f75=18, 4\", \"70\u00f77=10, 0\", \"69\u00f74=17, 1\", \"73\u00f76=12, 1\", \"56\u00f74=14, 0\"],\n];\n\nconst rowStep = 4; // rows with content are spaced 4 apart (1 content + 3 blank)\nfor (let blockIndex = 0; blockIndex < replacements.length; blockIndex++) {\n  const rowIndex = blockIndex * rowStep;\n  const rowValues = replacements[blockIndex];\n  for (let colIndex = 0; colIndex < rowValues.length; colIndex++) {\n    const cell = table.getCell(rowIndex, colIndex);\n    cell.value = rowValues[colIndex];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and every two-digit / one-digit division\n# practice problem in the worksheet table with the next day's\n# regenerated set.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date paragraph above the table -------------------\n$d.Paragraphs.Item(1).Range.Text = \"2024-09-16 Monday\"\n\n# --- 2. Update the practice-problem table ----------------------------\n$t = $d.Tables.Item(1)\n\n# Each 4-row \"block\" in the table holds one row of 5 problems followed\n# by 3 blank rows. Replacements are listed row-by-row, left-to-right,\n# in document order (1-based row/column indices, Word COM style).\n$replacements = @(\n    @(\"23\u00f74=5, 3\", \"82\u00f78=10, 2\", \"89\u00f78=11, 1\", \"82\u00f77=11, 5\", \"27\u00f73=9, 0\"),\n    @(\"71\u00f75=14, 1\", \"71\u00f78=8, 7\", \"38\u00f72=19, 0\", \"60\u00f77=8, 4\", \"91\u00f73=30, 1\"),\n    @(\"25\u00f73=8, 1\", \"79\u00f74=19, 3\", \"28\u00f76=4, 4\", \"80\u00f72=40, 0\", \"92\u00f79=10, 2\"),\n    @(\"16\u00f79=1, 7\", \"22\u00f78=2, 6\", \"62\u00f77=8, 6\", \"80\u00f77=11, 3\", \"60\u00f78=7, 4\"),\n    @(\"94\u00f75=18, 4\", \"70\u00f77=10, 0\", \"69\u00f74=17, 1\", \"73\u00f76=12, 1\", \"56\u00f74=14, 0\")\n)\n\n$rowStep = 4  # content rows are 4 apart (1 content row + 3 blank rows)\nfor ($blockIndex = 0; $blockIndex -lt $replacements.Count; $blockIndex++) {\n    $rowIndex = ($blockIndex * $rowStep) + 1\n    $rowValues = $replacements[$blockIndex]\n    for ($colIndex = 0; $colIndex -lt $rowValues.Count; $colIndex++) {\n        $cell = $t.Cell($rowIndex, $colIndex + 1)\n        $cell.Range.Text = $rowValues[$colIndex]\n    }\n}\n"}
